$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.195.87'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.680.32'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.93'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5271'
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2686'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06362'
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").Value = '1.690.96'
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.528'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5737'
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008231'
$ws.Range("E15").Value = '  -2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.44'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '26.226.55'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.868'
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.74'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.71'
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.227'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.007'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.15'
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.708'
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.85'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06411'
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.315'
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.569'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.562'
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.017'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6107'
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.744'
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01644'
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.172'
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("D40").Value = '1.094.00'
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8811'
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '1.834.70'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.48'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.088'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05268'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4280'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.014'
$ws.Range("E51").Value = '  -0.99%  '
